$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'24.927.27"
$ws.Range("E2").Value = "  +2.34%  "

$ws.Range("D3").Value = "'1.705.83"
$ws.Range("E3").Value = "  +1.58%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'316.07"
$ws.Range("E5").Value = "  +0.17%  "

$ws.Range("E6").Value = "  -0.04%  "

$ws.Range("D7").Value = "'0.3955"
$ws.Range("E7").Value = "  +2.02%  "

$ws.Range("D8").Value = "'0.4038"
$ws.Range("E8").Value = "  +1.16%  "

$ws.Range("D9").Value = "'1.485"
$ws.Range("E9").Value = "  +0.57%  "

$ws.Range("D10").Value = "'52.66"
$ws.Range("E10").Value = "  +0.95%  "

$ws.Range("D11").Value = "'1.002"
$ws.Range("E11").Value = "  -0.04%  "

$ws.Range("D12").Value = "'0.08810"
$ws.Range("E12").Value = "  +1.12%  "

$ws.Range("D13").Value = "'25.94"
$ws.Range("E13").Value = "  +0.43%  "

$ws.Range("D14").Value = "'7.466"
$ws.Range("E14").Value = "  +0.16%  "

$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "'0.00001358"
$ws.Range("E15").Value = "  +1.43%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'7.986"
$ws.Range("E16").Value = "  +0.29%  "

$ws.Range("D17").Value = "'1.712.97"
$ws.Range("E17").Value = "  +2.54%  "

$ws.Range("D18").Value = "'96.37"
$ws.Range("E18").Value = "  -1.28%  "

$ws.Range("D19").Value = "'0.07175"
$ws.Range("E19").Value = "  -0.28%  "

$ws.Range("D20").Value = "'20.64"
$ws.Range("E20").Value = "  +5.55%  "

$ws.Range("D21").Value = "'7.354"
$ws.Range("E21").Value = "  +1.67%  "

$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  -0.08%  "

$ws.Range("D23").Value = "'14.47"
$ws.Range("E23").Value = "  +2.56%  "

$ws.Range("D24").Value = "'24.939.13"
$ws.Range("E24").Value = "  +2.51%  "

$ws.Range("D25").Value = "'2.984"
$ws.Range("E25").Value = "  -0.38%  "

$ws.Range("D26").Value = "'2.350"
$ws.Range("E26").Value = "  +0.84%  "

$ws.Range("D27").Value = "'23.72"
$ws.Range("E27").Value = "  +5.80%  "

$ws.Range("D28").Value = "'6.164"
$ws.Range("E28").Value = "  +15.37%  "

$ws.Range("D29").Value = "'161.34"
$ws.Range("E29").Value = "  -3.37%  "

$ws.Range("D30").Value = "'150.57"
$ws.Range("E30").Value = "  +9.28%  "

$ws.Range("D31").Value = "'8.469"
$ws.Range("E31").Value = "  -1.49%  "

$ws.Range("D32").Value = "'2.554"
$ws.Range("E32").Value = "  +29.64%  "

$ws.Range("D33").Value = "'1.900.71"
$ws.Range("E33").Value = "  +2.49%  "

$ws.Range("D34").Value = "'0.08550"
$ws.Range("E34").Value = "  -2.16%  "

$ws.Range("D35").Value = "'0.03150"
$ws.Range("E35").Value = "  +7.29%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'1.047"
$ws.Range("E36").Value = "  +0.98%  "

$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").Value = "'7.188"
$ws.Range("E37").Value = "  -2.05%  "

$ws.Range("D38").Value = "'0.2856"
$ws.Range("E38").Value = "  +4.12%  "

$ws.Range("D39").Value = "'0.09540"
$ws.Range("E39").Value = "  +4.50%  "

$ws.Range("D40").Value = "'10.86"
$ws.Range("E40").Value = "  +1.07%  "

$ws.Range("D41").Value = "'0.8243"
$ws.Range("E41").Value = "  +4.68%  "

$ws.Range("D42").Value = "'14.00"
$ws.Range("E42").Value = "  +0.14%  "

$ws.Range("D43").Value = "'1.476"
$ws.Range("E43").Value = "  +0.55%  "

$ws.Range("D44").Value = "'17.45"
$ws.Range("E44").Value = "  +1.48%  "

$ws.Range("D45").Value = "'2.673"

$ws.Range("D46").Value = "'0.7384"
$ws.Range("E46").Value = "  +3.08%  "

$ws.Range("D47").Value = "'4.251"
$ws.Range("E47").Value = "  -0.01%  "

$ws.Range("E48").Value = "  -1.54%  "

$ws.Range("D49").Value = "'0.08740"
$ws.Range("E49").Value = "  +8.90%  "

$ws.Range("D50").Value = "'1.001"
$ws.Range("E50").Value = "  +0.09%  "

$ws.Range("D51").Value = "'139.07"
$ws.Range("E51").Value = "  -0.07%  "
